# 24.GFG - Swap two nodes in pair
# Adds a new row (row 20) to the Linked_List question tracker sheet:
#   A20 = "24/GFG"
#   B20 = "Swap Nodes in Pairs"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A20").Value = "24/GFG"
$ws.Range("B20").Value = "Swap Nodes in Pairs"

# Leave the newly-entered cell selected, matching the authoring session.
$ws.Range("B20").Select()
